$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: "Change in inventories"
$ws.Range("B6").Value = -224000000.0
$ws.Range("C6").Value = -240000000.0
$ws.Range("D6").Value = 30000000.0
$ws.Range("E6").Value = 13000000.0
$ws.Range("F6").Value = -36000000.0

# Row 7: "Change in payables and accrued liability"
$ws.Range("B7").Value = 412000000.0
$ws.Range("C7").Value = 237000000.0
$ws.Range("D7").Value = 431000000.0
$ws.Range("E7").Value = 29000000.0
$ws.Range("F7").Value = -93000000.0
